$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlRight horizontal alignment constant
$xlRight = -4152

# ---------------------------------------------------------------------
# Row 32 - new observation: Safer et al. (2010), eating, Emotional Eating
# Scale-Depression outcome
# ---------------------------------------------------------------------
$ws.Range("A32").Value = 15
$ws.Range("B32").Value = "Safer et al."
$ws.Range("B32").HorizontalAlignment = $xlRight
$ws.Range("C32").Value = 2010
$ws.Range("C32").HorizontalAlignment = $xlRight
$ws.Range("D32").Value = "eating"
$ws.Range("D32").Font.Name = "Times New Roman"
$ws.Range("D32").Font.Family = 1
$ws.Range("D32").HorizontalAlignment = $xlRight
$ws.Range("E32").Value = "Emotional Eating Scale-Depression"
$ws.Range("E32").HorizontalAlignment = $xlRight
$ws.Range("F32").Value = 43
$ws.Range("G32").Value = 43
$ws.Range("G32").HorizontalAlignment = $xlRight
$ws.Range("H32").Value = 0.26
$ws.Range("H32").HorizontalAlignment = $xlRight
$ws.Range("I32").Value = 0.05
$ws.Range("I32").HorizontalAlignment = $xlRight
$ws.Range("J32").Value = 0.22
$ws.Range("J32").HorizontalAlignment = $xlRight

# ---------------------------------------------------------------------
# Row 33 - new observation: Safer et al. (2010), eating, Emotional Eating
# Scale-Anxiety outcome
# ---------------------------------------------------------------------
$ws.Range("A33").Value = 15
$ws.Range("B33").Value = "Safer et al."
$ws.Range("B33").HorizontalAlignment = $xlRight
$ws.Range("C33").Value = 2010
$ws.Range("C33").HorizontalAlignment = $xlRight
$ws.Range("D33").Value = "eating"
$ws.Range("D33").Font.Name = "Times New Roman"
$ws.Range("D33").Font.Family = 1
$ws.Range("D33").HorizontalAlignment = $xlRight
$ws.Range("E33").Value = "Emotional Eating Scale-Anxiety"
$ws.Range("E33").HorizontalAlignment = $xlRight
$ws.Range("F33").Value = 43
$ws.Range("G33").Value = 43
$ws.Range("H33").Value = -0.34
$ws.Range("H33").HorizontalAlignment = $xlRight
$ws.Range("I33").Value = 0.05
$ws.Range("I33").HorizontalAlignment = $xlRight
$ws.Range("J33").Value = 0.12
$ws.Range("J33").HorizontalAlignment = $xlRight

# ---------------------------------------------------------------------
# Scroll / selection housekeeping to match the author's final view state
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D33").Select()
